# Roll the balance-sheet forward by one fiscal year:
#   - drop the oldest reporting period (1396/12) and its publish date,
#   - shift every remaining period one column to the left (D<-E<-F<-G<-H),
#   - populate the newly-opened rightmost column (H) with the 1401/12 figures.
# This mirrors the "update database" part of the commit; the period / publish
# date header labels shift the same way as the data rows below them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Row 8: "12 ماهه منتهی به ..." period headers
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# Row 9: "تاریخ انتشار" publish dates
$ws.Range("D9").Value = "1399-02-11 (8)"
$ws.Range("E9").Value = "1400-02-08 (8)"
$ws.Range("F9").Value = "1401-02-06 (9)"
$ws.Range("G9").Value = "1402-02-10 (8)"
$ws.Range("H9").Value = "1402-02-10 (2)"

# Row 12: موجودی نقد
$ws.Range("D12").Value = 25865
$ws.Range("E12").Value = 62627
$ws.Range("F12").Value = 262429
$ws.Range("G12").Value = 713802
$ws.Range("H12").Value = 615915

# Row 13: سرمایه گذاری کوتاه مدت
$ws.Range("D13").Value = 267
$ws.Range("E13").Value = 16
$ws.Range("F13").Value = 47
$ws.Range("G13").Value = 49981
$ws.Range("H13").Value = 191

# Row 14: دریافتنی‌های تجاری و سایر دریافتنی‌ها
$ws.Range("D14").Value = 207426
$ws.Range("E14").Value = 419488
$ws.Range("F14").Value = 441238
$ws.Range("G14").Value = 1389972
$ws.Range("H14").Value = 1753680

# Row 15: موجودی مواد و کالا
$ws.Range("D15").Value = 253741
$ws.Range("E15").Value = 276759
$ws.Range("F15").Value = 510116
$ws.Range("G15").Value = 871042
$ws.Range("H15").Value = 810766

# Row 16: پیش پرداخت ها
$ws.Range("D16").Value = 93388
$ws.Range("E16").Value = 83542
$ws.Range("F16").Value = 142817
$ws.Range("G16").Value = 312178
$ws.Range("H16").Value = 153396

# Row 17: دارایی های نگهداری شده برای فروش
$ws.Range("E17").Value = 629
$ws.Range("F17").Value = 0

# Row 18: جمع داراییهای جاری
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 843061
$ws.Range("F18").Value = 1356647
$ws.Range("G18").Value = 3336975
$ws.Range("H18").Value = 3333948

# Row 19: حسابها و اسناد دریافتنی تجاری بلند مدت
$ws.Range("D19").Value = 117
$ws.Range("E19").Value = 1406
$ws.Range("F19").Value = 715
$ws.Range("G19").Value = 1205
$ws.Range("H19").Value = 3940

# Row 20: سرمایه گذاریهای بلند مدت
$ws.Range("D20").Value = 0

# Row 21: سرمایه گذاری در املاک
$ws.Range("E21").Value = 43125

# Row 22: داراییهای ثابت مشهود
$ws.Range("D22").Value = 422500
$ws.Range("E22").Value = 491090
$ws.Range("F22").Value = 601734
$ws.Range("G22").Value = 1216060
$ws.Range("H22").Value = 2588976

# Row 23: داراییهای نامشهود
$ws.Range("D23").Value = 2721
$ws.Range("E23").Value = 2052
$ws.Range("F23").Value = 1936
$ws.Range("G23").Value = 1473
$ws.Range("H23").Value = 1291

# Row 25: سایر دارایی ها
$ws.Range("F25").Value = 14494
$ws.Range("G25").Value = 2000

# Row 26: جمع داراییهای غیرجاری
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 537673
$ws.Range("F26").Value = 662004
$ws.Range("G26").Value = 1263863
$ws.Range("H26").Value = 2639332

# Row 27: جمع داراییها
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 1380734
$ws.Range("F27").Value = 2018651
$ws.Range("G27").Value = 4600838
$ws.Range("H27").Value = 5973280

# Row 29: پرداختنی‌های تجاری و سایر پرداختنی‌ها
$ws.Range("D29").Value = 221809
$ws.Range("E29").Value = 276162
$ws.Range("F29").Value = 607689
$ws.Range("G29").Value = 1105379
$ws.Range("H29").Value = 1295990

# Row 31: پیش دریافتها
$ws.Range("D31").Value = 24340
$ws.Range("E31").Value = 3253
$ws.Range("F31").Value = 65303
$ws.Range("G31").Value = 113750
$ws.Range("H31").Value = 287883

# Row 32: ذخیره مالیات بر درامد
$ws.Range("D32").Value = 69976
$ws.Range("E32").Value = 133501
$ws.Range("F32").Value = 159734
$ws.Range("G32").Value = 323323
$ws.Range("H32").Value = 270919

# Row 33: سود سهام پیشنهادی و پرداختنی
$ws.Range("D33").Value = 5222
$ws.Range("E33").Value = 6945
$ws.Range("F33").Value = 33374
$ws.Range("G33").Value = 60408
$ws.Range("H33").Value = 53433

# Row 34: حصه جاری تسهیلات مالی دریافتی
$ws.Range("D34").Value = 48161
$ws.Range("E34").Value = 40385
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 600505
$ws.Range("H34").Value = 692130

# Row 37: جمع بدهیهای جاری
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = 460246
$ws.Range("F37").Value = 866100
$ws.Range("G37").Value = 2203365
$ws.Range("H37").Value = 2600355

# Row 39: حسابها و اسناد پرداختنی بلند مدت
$ws.Range("D39").Value = "-"

# Row 40: پیش دریافتهای غیرجاری
$ws.Range("D40").Value = 37797
$ws.Range("E40").Value = 0

# Row 41: تسهیلات مالی دریافتی بلند مدت
$ws.Range("D41").Value = 30713
$ws.Range("E41").Value = 43964
$ws.Range("F41").Value = 54202
$ws.Range("G41").Value = 84941
$ws.Range("H41").Value = 124529

# Row 42: جمع بدهیهای غیر جاری
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 43964
$ws.Range("F42").Value = 54202
$ws.Range("G42").Value = 84941
$ws.Range("H42").Value = 124529

# Row 43: جمع بدهیهای جاری و غیر جاری
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 504210
$ws.Range("F43").Value = 920302
$ws.Range("G43").Value = 2288306
$ws.Range("H43").Value = 2724884

# Row 45: سرمایه
$ws.Range("D45").Value = 474522
$ws.Range("H45").Value = 950000

# Row 47: وجوه دریافتی بابت افزایش سرمایه
$ws.Range("G47").Value = 437490
$ws.Range("H47").Value = 0

# Row 48: سهام خزانه
$ws.Range("F48").Value = -112365
$ws.Range("G48").Value = -89472
$ws.Range("H48").Value = -54800

# Row 49: صرف سهام خزانه
$ws.Range("D49").Value = 0

# Row 50: اندوخته قانونی
$ws.Range("D50").Value = 11054
$ws.Range("E50").Value = 28560
$ws.Range("F50").Value = 47452
$ws.Range("H50").Value = 95000

# Row 52: مازاد تجدید ارزیابی دارایی های غیر جاری نگه داری شده برای فروش
$ws.Range("D52").Value = "-"

# Row 53: مازاد تجدید ارزیابی داراییها
$ws.Range("H53").Value = 356000

# Row 54: اندوخته تسعیر ارز داراییها و بدهیهای شرکت های دولتی
$ws.Range("D54").Value = "-"

# Row 56: سود (زیان) انباشته
$ws.Range("D56").Value = 82431
$ws.Range("E56").Value = 373442
$ws.Range("F56").Value = 688740
$ws.Range("G56").Value = 1442540
$ws.Range("H56").Value = 1902196

# Row 57: جمع حقوق صاحبان سهام
$ws.Range("D57").Value = 0
$ws.Range("E57").Value = 876524
$ws.Range("F57").Value = 1098349
$ws.Range("G57").Value = 2312532
$ws.Range("H57").Value = 3248396

# Row 58: جمع بدهیها و حقوق صاحبان سهام
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 1380734
$ws.Range("F58").Value = 2018651
$ws.Range("G58").Value = 4600838
$ws.Range("H58").Value = 5973280
